$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 3 for columns D, J, K, L, M, P

$d2 = $ws.Range("D2").Value2
$j2 = $ws.Range("J2").Value2
$k2 = $ws.Range("K2").Value2
$l2 = $ws.Range("L2").Value2
$m2 = $ws.Range("M2").Value2
$p2 = $ws.Range("P2").Value2

$d3 = $ws.Range("D3").Value2
$j3 = $ws.Range("J3").Value2
$k3 = $ws.Range("K3").Value2
$l3 = $ws.Range("L3").Value2
$m3 = $ws.Range("M3").Value2
$p3 = $ws.Range("P3").Value2

$ws.Range("D2").Value = $d3
$ws.Range("J2").Value = $j3
$ws.Range("K2").Value = $k3
$ws.Range("L2").Value = $l3
$ws.Range("M2").Value = $m3
$ws.Range("P2").Value = $p3

$ws.Range("D3").Value = $d2
$ws.Range("J3").Value = $j2
$ws.Range("K3").Value = $k2
$ws.Range("L3").Value = $l2
$ws.Range("M3").Value = $m2
$ws.Range("P3").Value = $p2

$wb.Save()
